$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Row 7: Credenciales_ChatGPT / ProyectoRPA / Credenciales de Acceso para ChatGPT
$ws.Range("A7").Value = "Credenciales_ChatGPT"
$ws.Range("B7").Value = "ProyectoRPA"
$ws.Range("B7").Style = "Hipervínculo"
$ws.Range("C7").Value = "Credenciales de Acceso para ChatGPT"

# Row 8: DireccionURL / https://chatgpt.com/auth/login (hyperlink) / Dirección de ChatGPT
$ws.Range("A8").Value = "DireccionURL"
$ws.Range("B8").Value = "https://chatgpt.com/auth/login"
$ws.Hyperlinks.Add($ws.Range("B8"), "https://chatgpt.com/auth/login")
$ws.Range("B8").Style = "Hipervínculo"
$ws.Range("C8").Value = "Dirección de ChatGPT"
$ws.Range("C8").WrapText = $true

# Row 9: DireccionURLIdealista / https://www.idealista.com/ (hyperlink) / Dirección de Idealista
$ws.Range("A9").Value = "DireccionURLIdealista"
$ws.Range("B9").Value = "https://www.idealista.com/"
$ws.Hyperlinks.Add($ws.Range("B9"), "https://www.idealista.com/")
$ws.Range("B9").Style = "Hipervínculo"
$ws.Range("C9").Value = "Dirección de Idealista"

# Row 10: empty cell C10 keeps the wrapped/description style
$ws.Range("C10").WrapText = $true

# Last clicked/selected cell in the Settings sheet
$ws.Range("A6").Select()
